$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $text)
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '245.97'
Set-TextValue $ws 'D4' '5.280'
Set-TextValue $ws 'D5' '0.05802'
Set-TextValue $ws 'D7' '3.132'
Set-TextValue $ws 'D8' '0.8167'
Set-TextValue $ws 'D9' '0.8526'
Set-TextValue $ws 'B10' 'One'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D10' '0.009791'
Set-TextValue $ws 'E10' '9OneONEBestin24h'
Set-TextValue $ws 'B11' 'WazirX'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D11' '0.1360'
Set-TextValue $ws 'E11' '10WazirXWRX'
Set-TextValue $ws 'B12' 'MandalaExchangeToken'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D12' '0.06938'
Set-TextValue $ws 'E12' '11MandalaExchangeTokenMDX'
Set-TextValue $ws 'B13' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D13' '0.03134'
Set-TextValue $ws 'E13' '12LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws 'B14' 'BitrueCoin'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D14' '0.02868'
Set-TextValue $ws 'E14' '13BitrueCoinBTR'
Set-TextValue $ws 'B15' 'BitMartToken'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D15' '0.09399'
Set-TextValue $ws 'E15' '14BitMartTokenBMX'
Set-TextValue $ws 'B16' 'MCDex'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws 'D16' '3.743'
Set-TextValue $ws 'E16' '15MCDexMCB'
Set-TextValue $ws 'B17' 'BitForexToken'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D17' '0.001525'
Set-TextValue $ws 'E17' '16BitForexTokenBF'
Set-TextValue $ws 'B18' 'CoinExToken'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws 'D18' '0.04669'
Set-TextValue $ws 'E18' '17CoinExTokenCET'
Set-TextValue $ws 'D19' '0.006262'
Set-TextValue $ws 'D20' '0.001237'
Set-TextValue $ws 'D21' '0.004627'
Set-TextValue $ws 'D23' '3.501'
Set-TextValue $ws 'D28' '0.0002329'
Set-TextValue $ws 'D40' '0.03660'
Set-TextValue $ws 'B41' 'KickToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D41' '0.006264'
Set-TextValue $ws 'E41' '40KickTokenKICK'
Set-TextValue $ws 'B42' 'BKEXToken'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D42' '0.1057'
Set-TextValue $ws 'E42' '41BKEXTokenBKK'
Set-TextValue $ws 'B43' 'CEJI'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws 'D43' '0.003397'
Set-TextValue $ws 'E43' '42CEJICEJI'
Set-TextValue $ws 'D44' '0.007468'
Set-TextValue $ws 'D45' '0.00005271'
Set-TextValue $ws 'D46' '0.00000000749'
Set-TextValue $ws 'D47' '0.3697'
Set-TextValue $ws 'E47' '46CoinbaseStockTokenCOINWorstin24h'
Set-TextValue $ws 'D48' '0.002197'
Set-TextValue $ws 'D49' '0.00002098'
Set-TextValue $ws 'D50' '0.0001998'
